$d = $word.ActiveDocument

# 1. Move the "_GoBack" bookmark out of the end of the "HTML was initially..."
#    paragraph into the middle of the title ("True or False Ques|tions"),
#    splitting the title run in two. Word treats "_GoBack" as a singleton
#    bookmark, so re-adding it here removes the old occurrence automatically.
$titleSplit = $d.Content
$titleSplit.Find.Execute("True or False Ques", $true, $false, $false, $false, $false, `
                          $true, 1, $false, "", 0) | Out-Null
$titleSplit.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $titleSplit) | Out-Null

# 2. Trim the trailing " F" (a stray space run followed by an
#    underlined "F" run) from the end of the "WWW was used..." paragraph.
$d.Content.Find.Execute(
    "PRESENCE. F",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "PRESENCE.",
    2
) | Out-Null
